$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (26) down to the new rows (27-31)
# so that column A keeps its bold/border style and column E keeps its date number format.
$ws.Range("A26:V26").Copy() | Out-Null
$ws.Range("A27:V31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 27
$ws.Range("A27").Value2 = 26
$ws.Range("B27").Value = "kuwait"
$ws.Range("C27").Value = "premier-league"
$ws.Range("D27").Value = "2023-2024"
$ws.Range("E27").Value2 = 45226.75
$ws.Range("F27").Value = "Al Shabab"
$ws.Range("G27").Value2 = 1
$ws.Range("H27").Value = "Al Salmiya"
$ws.Range("I27").Value2 = 2
$ws.Range("J27").Value2 = 2.61
$ws.Range("K27").Value = "26/10/2023 15:43"
$ws.Range("L27").Value2 = 2.97
$ws.Range("M27").Value = "27/10/2023 00:55"
$ws.Range("N27").Value2 = 3.28
$ws.Range("O27").Value = "26/10/2023 15:43"
$ws.Range("P27").Value2 = 3.68
$ws.Range("Q27").Value = "27/10/2023 16:06"
$ws.Range("R27").Value2 = 2.28
$ws.Range("S27").Value = "26/10/2023 15:43"
$ws.Range("T27").Value2 = 2.08
$ws.Range("U27").Value = "27/10/2023 00:55"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-shabab-al-salmiya/jevNQ9mP/"

# Row 28
$ws.Range("A28").Value2 = 27
$ws.Range("B28").Value = "kuwait"
$ws.Range("C28").Value = "premier-league"
$ws.Range("D28").Value = "2023-2024"
$ws.Range("E28").Value2 = 45227.6875
$ws.Range("F28").Value = "Al Jahra"
$ws.Range("G28").Value2 = 2
$ws.Range("H28").Value = "Al Arabi"
$ws.Range("I28").Value2 = 3
$ws.Range("J28").Value2 = 4.52
$ws.Range("K28").Value = "27/10/2023 15:43"
$ws.Range("L28").Value2 = 6.11
$ws.Range("M28").Value = "28/10/2023 01:48"
$ws.Range("N28").Value2 = 3.87
$ws.Range("O28").Value = "27/10/2023 15:43"
$ws.Range("P28").Value2 = 4.15
$ws.Range("Q28").Value = "28/10/2023 14:32"
$ws.Range("R28").Value2 = 1.54
$ws.Range("S28").Value = "27/10/2023 15:43"
$ws.Range("T28").Value2 = 1.46
$ws.Range("U28").Value = "28/10/2023 01:48"
$ws.Range("V28").Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-jahra-al-arabi-kuwait/AuLgLmBn/"

# Row 29
$ws.Range("A29").Value2 = 28
$ws.Range("B29").Value = "kuwait"
$ws.Range("C29").Value = "premier-league"
$ws.Range("D29").Value = "2023-2024"
$ws.Range("E29").Value2 = 45227.80555555555
$ws.Range("F29").Value = "Al-Fahaheel"
$ws.Range("G29").Value2 = 4
$ws.Range("H29").Value = "Al Naser"
$ws.Range("I29").Value2 = 2
$ws.Range("J29").Value2 = 2.68
$ws.Range("K29").Value = "27/10/2023 15:43"
$ws.Range("L29").Value2 = 3.63
$ws.Range("M29").Value = "28/10/2023 18:59"
$ws.Range("N29").Value2 = 3.29
$ws.Range("O29").Value = "27/10/2023 15:43"
$ws.Range("P29").Value2 = 3.52
$ws.Range("Q29").Value = "28/10/2023 18:59"
$ws.Range("R29").Value2 = 2.22
$ws.Range("S29").Value = "27/10/2023 15:43"
$ws.Range("T29").Value2 = 1.9
$ws.Range("U29").Value = "28/10/2023 18:59"
$ws.Range("V29").Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-fahaheel-al-naser/h4QlMTet/"

# Row 30
$ws.Range("A30").Value2 = 29
$ws.Range("B30").Value = "kuwait"
$ws.Range("C30").Value = "premier-league"
$ws.Range("D30").Value = "2023-2024"
$ws.Range("E30").Value2 = 45228.64583333334
$ws.Range("F30").Value = "Khaitan"
$ws.Range("G30").Value2 = 0
$ws.Range("H30").Value = "Kazma SC"
$ws.Range("I30").Value2 = 3
$ws.Range("J30").Value2 = 3.8
$ws.Range("K30").Value = "28/10/2023 09:43"
$ws.Range("L30").Value2 = 4.06
$ws.Range("M30").Value = "29/10/2023 14:54"
$ws.Range("N30").Value2 = 3.65
$ws.Range("O30").Value = "28/10/2023 09:43"
$ws.Range("P30").Value2 = 3.98
$ws.Range("Q30").Value = "29/10/2023 14:54"
$ws.Range("R30").Value2 = 1.69
$ws.Range("S30").Value = "28/10/2023 09:43"
$ws.Range("T30").Value2 = 1.7
$ws.Range("U30").Value = "28/10/2023 17:54"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/kuwait/premier-league/khaitan-kazma-sc/xhju2QAO/"

# Row 31
$ws.Range("A31").Value2 = 30
$ws.Range("B31").Value = "kuwait"
$ws.Range("C31").Value = "premier-league"
$ws.Range("D31").Value = "2023-2024"
$ws.Range("E31").Value2 = 45228.76388888889
$ws.Range("F31").Value = "Al Qadisiya"
$ws.Range("G31").Value2 = 1
$ws.Range("H31").Value = "Al Kuwait"
$ws.Range("I31").Value2 = 1
$ws.Range("J31").Value2 = 3.63
$ws.Range("K31").Value = "28/10/2023 09:43"
$ws.Range("L31").Value2 = 3.59
$ws.Range("M31").Value = "29/10/2023 04:33"
$ws.Range("N31").Value2 = 3.6
$ws.Range("O31").Value = "28/10/2023 09:43"
$ws.Range("P31").Value2 = 3.69
$ws.Range("Q31").Value = "29/10/2023 16:23"
$ws.Range("R31").Value2 = 1.77
$ws.Range("S31").Value = "28/10/2023 09:43"
$ws.Range("T31").Value2 = 1.85
$ws.Range("U31").Value = "29/10/2023 04:33"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-qadisiya-al-kuwait/lWJcK7Qh/"

$ws.Range("A1").Select() | Out-Null
Write-Output $ws.UsedRange.Address()
